$d = $word.ActiveDocument

# Map of distinctive paragraph-text prefixes -> full replacement text.
# Using $p.Range.Text = "..." (rather than Find.Execute replacement) keeps
# the run's existing formatting / xml:space="preserve" intact in this engine.
$edits = @(
    @{ Match = "Ah, well. At least the day*"; New = "Well, at least the day’s over. I get up and stretch, ready to go home and relax…" },
    @{ Match = "There it is. Well, might as well get it over with.*"; New = "There it is. Might as well get it over with, I guess." },
    @{ Match = "Asher (neutral curious): Sure, that’d be great. Your mom will be okay with it?*"; New = "Asher (neutral curious): Sure, that’d be great. Will your mom be okay with it?" },
    @{ Match = "Petra (neutral raised_eyebrow): You sure?*"; New = "Petra: You sure?" }
)

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    foreach ($edit in $edits) {
        if ($t -like $edit.Match) {
            $p.Range.Text = $edit.New
        }
    }
}
